# Update cryptos price (D) and volume-change (E) columns per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column price values that are plain text and never look like numbers
# to Excel (contain two or more '.' separators) -- safe to assign directly.
$textPriceUpdates = @{
    'D2' = '67.604.80'
    'D3' = '2.681.08'
    'D9' = '2.680.40'
    'D15' = '3.169.65'
    'D17' = '67.595.08'
    'D18' = '2.667.39'
    'D28' = '2.806.08'
}
foreach ($ref in $textPriceUpdates.Keys) {
    $ws.Range($ref).Value = $textPriceUpdates[$ref]
}

# D-column price values that Excel would otherwise auto-convert to a
# number (single decimal point) -- force them to stay text, matching the
# original inline-string cell type, then restore the default cell style
# so no new number-format style is left applied to the cell.
$numericLookingPriceUpdates = @{
    'D4' = '1.00'
    'D5' = '600.37'
    'D6' = '165.70'
    'D10' = '0.143'
    'D14' = '27.83'
    'D19' = '11.77'
    'D20' = '7.65'
    'D21' = '364.27'
    'D23' = '4.82'
    'D26' = '71.19'
    'D27' = '10.11'
    'D29' = '0.0000102'
    'D31' = '558.21'
    'D32' = '8.00'
    'D38' = '19.57'
    'D39' = '154.77'
    'D46' = '40.40'
    'D48' = '0.591'
    'D49' = '153.48'
    'D51' = '1.72'
}
foreach ($ref in $numericLookingPriceUpdates.Keys) {
    $ws.Range($ref).NumberFormat = '@'
}
foreach ($ref in $numericLookingPriceUpdates.Keys) {
    $ws.Range($ref).Value = $numericLookingPriceUpdates[$ref]
}
foreach ($ref in $numericLookingPriceUpdates.Keys) {
    $ws.Range($ref).Style = 'Normal'
}

# E-column Volume(1h) percentage-change values (always plain text,
# never numeric-looking because of the leading/trailing spaces and '%').
$volumeUpdates = @{
    'E2' = '  -1.27%  '
    'E3' = '  -0.35%  '
    'E4' = '  -0.04%  '
    'E5' = '  +0.33%  '
    'E6' = '  +3.67%  '
    'E7' = '  +0.02%  '
    'E8' = '  +0.70%  '
    'E9' = '  -0.36%  '
    'E10' = '  +1.65%  '
    'E11' = '  +1.07%  '
    'E12' = '  -0.22%  '
    'E13' = '  -1.22%  '
    'E14' = '  -1.33%  '
    'E15' = '  -0.30%  '
    'E16' = '  -2.16%  '
    'E17' = '  -1.18%  '
    'E18' = '  -0.57%  '
    'E19' = '  -0.64%  '
    'E20' = '  +0.45%  '
    'E21' = '  -0.72%  '
    'E22' = '  -3.27%  '
    'E23' = '  -0.95%  '
    'E24' = '  -3.70%  '
    'E25' = '  +0.10%  '
    'E26' = '  -4.42%  '
    'E27' = '  +1.14%  '
    'E29' = '  -1.86%  '
    'E30' = '  +0.00%  '
    'E31' = '  -2.31%  '
    'E32' = '  -2.62%  '
    'E34' = '  +0.06%  '
    'E35' = '  -1.55%  '
    'E36' = '  +0.03%  '
    'E37' = '  -5.17%  '
    'E38' = '  -1.38%  '
    'E39' = '  -4.35%  '
    'E40' = '  -1.04%  '
    'E41' = '  -1.51%  '
    'E42' = '  -4.02%  '
    'E43' = '  +0.53%  '
    'E44' = '  -4.73%  '
    'E46' = '  -0.35%  '
    'E47' = '  -5.35%  '
    'E48' = '  -0.97%  '
    'E49' = '  -2.41%  '
    'E50' = '  -3.86%  '
    'E51' = '  -2.56%  '
}
foreach ($ref in $volumeUpdates.Keys) {
    $ws.Range($ref).Value = $volumeUpdates[$ref]
}
